# Generate Report for Handoff
# Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the 482a09c7-74a0-458d-b243-195968d1d230.md entry across
# the Overview, zh-cn and de-de sheets (row 6 in each table).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G6").Value = "2016-08-27 06:39:10"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H6").Value = "2016-08-27 06:39:06"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H6").Value = "2016-08-27 06:39:10"
